$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.2544960331724724; C = 0.047201678941164; D = 0.07895109456468674; E = 0.1657639099158317; G = 0.4575342760136678; H = 0.6081937496374437; I = 0.4751391974968548; K = 0.2682970138353369; M = 0.2191956777404087; N = 1.328724512029542; O = 2.087546030305489 }
    3  = @{ B = 0.2232889452707525; C = 0.04226020517934614; D = 0.07156561347113666; E = 0.1547355734945413; G = 0.4569457636006646; H = 0.6113509066354439; I = 0.4785726943279833; K = 0.2340793334095821; M = 0.1967982620147097; N = 1.343306159277571; O = 2.092659347775836 }
    4  = @{ B = 0.2041033457231549; C = 0.03920596939185828; D = 0.06706312771380851; E = 0.1480777763533609; G = 0.4568957363290878; H = 0.6135501414614097; I = 0.4809330360301907; K = 0.2130136801549014; M = 0.1831123702986375; N = 1.352717766501861; O = 2.096995970917874 }
    5  = @{ B = 0.1962793947235184; C = 0.03795631765034102; D = 0.0652364718432068;  E = 0.1453930639361189; G = 0.4569535773651552; H = 0.6145119594796;    I = 0.4819582939928821; K = 0.2044156547175646; M = 0.177551919679523;  N = 1.356668324518956; O = 2.0990641596957 }
    6  = @{ B = 0.1949799048838088; C = 0.03774851159772652; D = 0.0649336500352149;  E = 0.144948979056025;  G = 0.4569679049164677; H = 0.6146756329038254; I = 0.4821323662349499; K = 0.2029871508151189; M = 0.1766296177677873; N = 1.357331272006506; O = 2.099425758679246 }
    7  = @{ B = 0.2039978515155099; C = 0.03918913644696431; D = 0.06703845976866774; E = 0.1480414546588875; G = 0.4568961997289946; H = 0.6135628471498151; I = 0.4809466063266399; K = 0.2128977786553321; M = 0.1830373126396694; N = 1.352770578487244; O = 2.097022644611471 }
    8  = @{ B = 0.2437411391810258; C = 0.0455020644685078;  D = 0.07639791282230135; E = 0.1619376385490625; G = 0.4572667063515539; H = 0.6092282527562105; I = 0.4762707259472911; K = 0.2565106083622481; M = 0.2114593067073116; N = 1.333657108271116; O = 2.089060641599886 }
    9  = @{ B = 0.3214684471953433; C = 0.05772055025957457; D = 0.09500672561759416; E = 0.1900995018486356; G = 0.4604668039515047; H = 0.6027949419234773; I = 0.4691026872680943; K = 0.3415764362759148; M = 0.2677225128029548; N = 1.299813348911132; O = 2.082948326736158 }
    10 = @{ B = 0.378430895990931;  C = 0.06659835895675315; D = 0.1088344158860366;  E = 0.2113613485988068; G = 0.4643315831875583; H = 0.5993261354172859; I = 0.4650574648527694; K = 0.4037795489812765; M = 0.3093883727763327; N = 1.277165647630588; O = 2.084257620356368 }
    11 = @{ B = 0.4043102250342656; C = 0.07061544525657837; D = 0.1151589713017245;  E = 0.2211614559501811; G = 0.4664197482418757; H = 0.5980207540376909; I = 0.4634825997908294; K = 0.4320105026572492; M = 0.3284166057665914; N = 1.267344026769358; O = 2.086114584522761 }
    12 = @{ B = 0.4141048920441222; C = 0.07213349224194587; D = 0.1175588212748266;  E = 0.2248911376040326; G = 0.4672580265950472; H = 0.5975655992791076; I = 0.4629244117077143; K = 0.4426910253042138; M = 0.3356328265992943; N = 1.263694021831136; O = 2.086999248842801 }
    13 = @{ B = 0.4119956777907134; C = 0.07180669392850803; D = 0.1170417544925471;  E = 0.2240870539582502; G = 0.4670753730306245; H = 0.5976618836746752; I = 0.463042929066809;  K = 0.4403912331191293; M = 0.3340782115369052; N = 1.264477035515199; O = 2.086800648098574 }
    14 = @{ B = 0.4051161467475026; C = 0.07074039903189089; D = 0.1153563109274529;  E = 0.2214679253763947; G = 0.4664877608034459; H = 0.5979825234546752; I = 0.4634359120089862; K = 0.4328893974979735; M = 0.329010075500868;  N = 1.267042350777259; O = 2.08618372890075 }
    15 = @{ B = 0.4009015336845323; C = 0.07008685272347748; D = 0.1143245618548434;  E = 0.2198660605529028; G = 0.4661340240554637; H = 0.5981840239868745; I = 0.4636815982129256; K = 0.4282929989529691; M = 0.3259070807847024; N = 1.268622698385693; O = 2.085829483787649 }
    16 = @{ B = 0.3767389104439474; C = 0.06633539666854915; D = 0.1084217766704541;  E = 0.2107234801539803; G = 0.4642017646087169; H = 0.5994169275664092; I = 0.4651657272212617; K = 0.4019332278957961; M = 0.3081463239271329; N = 1.277817197360243; O = 2.084161652909756 }
    17 = @{ B = 0.3619070692306252; C = 0.06402846549823948; D = 0.1048093512782913;  E = 0.2051477185624933; G = 0.463100974359449;  H = 0.6002430680058097; I = 0.4661441670145265; K = 0.3857452163320261; M = 0.2972696635412078; N = 1.283580953599534; O = 2.083461621814081 }
    18 = @{ B = 0.3533730939654163; C = 0.06269956143408706; D = 0.1027348084086697;  E = 0.2019527233928713; G = 0.462498891145998;  H = 0.6007439021940257; I = 0.4667319120098696; K = 0.3764281598782304; M = 0.2910206844657779; N = 1.286941376867492; O = 2.083177702273616 }
    19 = @{ B = 0.350483120754177;  C = 0.06224927219078324; D = 0.1020329598691916;  E = 0.2008730142418926; G = 0.4623003688255949; H = 0.6009178842300145; I = 0.4669352006082264; K = 0.3732725252750129; M = 0.2889060887062342; N = 1.28808693026945;  O = 2.083101959477119 }
    20 = @{ B = 0.3634862677493231; C = 0.06427425149615829; D = 0.105193566331252;   E = 0.2057400209089622; G = 0.4632149400005261; H = 0.6001524684732402; I = 0.4660374257032132; K = 0.3874690969247752; M = 0.2984267795690414; N = 1.28296270682658;  O = 2.083523853602316 }
    21 = @{ B = 0.4071369788999846; C = 0.07105368119940181; D = 0.115851234570286;   E = 0.2222367207955784; G = 0.466659066073305;  H = 0.5978872811273561; I = 0.4633194470195789; K = 0.4350931418657069; M = 0.330498421344565;  N = 1.266286976038415; O = 2.086360007203638 }
    22 = @{ B = 0.4356342720262774; C = 0.07546612642819639; D = 0.1228450463662654;  E = 0.2331267335767961; G = 0.4691871088775343; H = 0.596635117233518;  I = 0.4617656398134073; K = 0.4661601148239072; M = 0.3515211073252544; N = 1.255791980973057; O = 2.089271398011505 }
    23 = @{ B = 0.4204277302989965; C = 0.07311281224362176; D = 0.1191097346244447;  E = 0.2273045417459088; G = 0.4678124653669045; H = 0.5972825454468165; I = 0.4625745643259762; K = 0.4495845675910175; M = 0.3402952405684161; N = 1.261356407201287; O = 2.087620715603975 }
    24 = @{ B = 0.362772333716606;  C = 0.06416313979961785; D = 0.105019855575037;   E = 0.2054722080651672; G = 0.4631633202751004; H = 0.6001933479192019; I = 0.4660856048712851; K = 0.3866897627198398; M = 0.2979036344619033; N = 1.283242070531576; O = 2.083495349378381 }
    25 = @{ B = 0.3004652270059012; C = 0.05443246812561142; D = 0.08994521285126211; E = 0.1823817702852324; G = 0.4593356889427014; H = 0.6043142888198787; I = 0.4708274640905792; K = 0.3186145489814578; M = 0.2524445402382227; N = 1.308579585464428; O = 2.083583892296332 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
